$d = $word.ActiveDocument

# 1. Update activation date: 01/01/2018 -> 01/01/2025
$ok1 = $d.Content.Find.Execute("Ativação: 01/01/2018", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2025", 2)
Write-Host "Step1 (activation date) replaced: $ok1"

# 2. Fill in the empty English "Objetivos" paragraph (italic run currently has no text)
$objParaEn = $d.Paragraphs(7)
$objParaEn.Range.Text = "Provide the student with theoretical knowledge in the technological processes of preparing fermented and distilled beverages. Knowledge of raw materials, wort preparation, fermentation technology, distillation, beverage maturation."
Write-Host "Step2 (English objectives) set."

# 3. Update "Programa resumido" Portuguese paragraph
$ok3 = $d.Content.Find.Execute("Generalidades, tipos de bebidas e preparação do mosto; tipos de leveduras; fermentação alcoólica. Acabamento: maturação, filtração, destilação, envelhecimento. Análise química e sensorial.", $true, $false, $false, $false, $false, $true, 1, $false, "Generalidades, tipos de bebidas alcoólicas e preparação do mosto; tipos de leveduras; fermentação alcoólica. Acabamento: maturação, filtração, destilação, envelhecimento. Bebidas fermentadas não-alcoólicas: generalidades e processo produtivo. Análise química e sensorial.", 2)
Write-Host "Step3 (resumo PT) replaced: $ok3"

# 4. Update "Programa resumido" English paragraph
$ok4 = $d.Content.Find.Execute("Generalities, types of beverages and wort preparation; types of yeasts; alcoholic fermentation. Finishing: maturation, filtration, distillation, aging. Chemical and sensorial analyzes.", $true, $false, $false, $false, $false, $true, 1, $false, "General information, types of alcoholic beverages and must preparation; types of yeast; alcoholic fermentation. Finishing: maturation, filtration, distillation, aging. Non-alcoholic fermented beverages: generalities and production process. Chemical and sensory analysis of beverages.", 2)
Write-Host "Step4 (resumo EN) replaced: $ok4"

# 5. Update "Programa" Portuguese paragraph (the source run mixes in a Greek question
#    mark confusable (U+037E) for ';' in several spots — safer to replace the whole
#    paragraph range directly by index instead of Find/Replace text-matching it)
$programaPt = $d.Paragraphs(13)
Write-Host "Step5 target (old PT programa) starts: $($programaPt.Range.Text.Substring(0, [Math]::Min(40, $programaPt.Range.Text.Length)))"
$programaPt.Range.Text = "1. Generalidades: origem das bebidas; matérias-primas. 2. Tipos de bebidas: bebidas fermentadas (cervejas, fermentados de frutas, vinhos, sidras), bebidas destiladas (aguardentes, destilados de vinhos, graspa, pisco, rum, tequila, tiquira, uísque), bebidas retificadas (vodka, gim), bebidas obtidas por misturas (licores, sangria, cooler). 3. Preparação do mosto: pé de cuba, características físicoquímicas, correção do mosto. 4. Tipos de leveduras: leveduras selvagens, leveduras mistas, leveduras selecionadas. 5. Fermentação alcoólica: controle da fermentação, rendimento da fermentação, produtos secundários. 6. Acabamento: controle da maturação, destilação em alambiques e em colunas, determinação do grau alcoólico, armazenamento, tipos de madeiras, cor, volume, composição da bebida, legislação. 7. Preparação de bebiidas fermentadas não alcoólicas (iogurte, leite fermentado, kefir e outras). 8.Análise química: composição da bebida, legislação. 8. Análise sensorial: aromas das bebidas e aceitação."

# 6. Update "Programa" English paragraph (same confusable-character issue)
$programaEn = $d.Paragraphs(14)
Write-Host "Step6 target (old EN programa) starts: $($programaEn.Range.Text.Substring(0, [Math]::Min(40, $programaEn.Range.Text.Length)))"
$programaEn.Range.Text = "1. General: origin of the beverages; 2. Types of beverages: fermented beverages (beers, fermented fruits, wines, ciders), distilled beverages (cachaça, wine distillates, graspa, pisco, rum, tequila, tiquira, whiskey), rectified beverages (vodka, gin), beverages obtained by mixing (liqueurs, sangria, cooler). 3. Preparation of the must: foot of vat, physicochemical characteristics, correction of the wort. 4. Types of yeast: wild yeast, mixed yeast, selected yeast. 5. Alcoholic fermentation: fermentation control, fermentation yield, secondary products. 6. Finishing: control of maturation, distillation in stills and columns, determination of alcoholic content, storage, types of wood, color, volume, composition of the drink, legislation. 7. Preparation of non-alcoholic fermented beverages (yogurt, fermented milk, kefir and others). 8. Chemical analysis: composition of the beverages, legislation. 9. Sensory analysis: beverage aromas and acceptance"

# 7. Update Bibliography paragraph (same confusable-character issue)
$biblio = $d.Paragraphs(18)
Write-Host "Step7 target (old biblio) starts: $($biblio.Range.Text.Substring(0, [Math]::Min(40, $biblio.Range.Text.Length)))"
$biblio.Range.Text = "1) DA SILVA, N., JUNQUEIRA, V. C. A., DE ARRUDA SILVEIRA, N. F., TANIWAKI, M. H., GOMES, R. A. R., OKAZAKI, M. M. Manual de métodos de análise microbiológica de alimentos e água. Editora Blucher, 2017. 2) DA-SILVA, R.; LAGO-VANZELA, E. S.; BAFFI, M. A. Uvas e vinhos: química, bioquímica e microbiologia. São Paulo, Editora Senac, 2015. 3) DE OLIVEIRA MORAES, I. Biotecnologia Industrial: biotecnologia na produção de alimentos. Vol. 4. 2ª Ed. Editora Blucher, 2021. 4) MARTIN, J. G. P., DE DEA LINDNER, J. Microbiologia de alimentos fermentados. Editora Blucher, 2022. 5) MENEZES e SILVA, C.H.P. Microbiologia da cerveja - Do básico ao avançado, o guia definitivo. Editora LF, 2019. 6) MUXEL, A. A. Química da Cerveja: Uma Abordagem Química e Bioquímica das Matérias-Primas, Processo de Produção e da Composição dos Compostos de Sabores da Cerveja. Editora Appris, 2022. 7) VENTURINI FILHO, W. G. Bebidas alcoólicas: ciência e tecnologia. Vol. 1. Editora Blucher, 2021."

Write-Host "All edits applied."
